# Actualización automática 2025-10-30 16:30:08
$wb = $excel.ActiveWorkbook

$wsVentasGrupo   = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsVentaMensual  = $wb.Worksheets.Item("VENTA MENSUAL")
$wsCumplimiento  = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# ---- Hoja "VENTAS POR GRUPO" ----
$wsVentasGrupo.Range("C4").Value = 2073.6
$wsVentasGrupo.Range("M4").Value = 5370.43
$wsVentasGrupo.Range("M17").Value = 1161.46
$wsVentasGrupo.Range("L29").Value = 2262.54
$wsVentasGrupo.Range("M37").Value = 4200.62
$wsVentasGrupo.Range("M56").Value = "15 de 54"

# ---- Hoja "VENTA MENSUAL" ----
$wsVentaMensual.Range("F4").Value = 8395.700000000001
$wsVentaMensual.Range("F17").Value = 1161.46
$wsVentaMensual.Range("F29").Value = 11218.51
$wsVentaMensual.Range("F37").Value = 6051.97
$wsVentaMensual.Range("F60").Value = 91195.97

# ---- Hoja "CUMPLIMIENTO MENSUAL" ----
$wsCumplimiento.Range("D2").Value = 2073.6
$wsCumplimiento.Range("E2").Value = 4123.984029436589
$wsCumplimiento.Range("F2").Value = 0.3345819903612516

$wsCumplimiento.Range("D11").Value = 11397.6
$wsCumplimiento.Range("E11").Value = 433.3999999999996
$wsCumplimiento.Range("F11").Value = 0.9633674245625898

$wsCumplimiento.Range("D12").Value = 52183.43
$wsCumplimiento.Range("E12").Value = 479.6900000000023
$wsCumplimiento.Range("F12").Value = 0.9908913486325914

$wsCumplimiento.Range("D14").Value = 85946.73999999999
$wsCumplimiento.Range("E14").Value = 13069.76661190614
$wsCumplimiento.Range("F14").Value = 0.8680041635569621
